$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns F (percentil_dist_50), G (percentil_densidad_25),
# K (densidad) and Z (timestamp) for each data row (2-112), derived from the
# committed change to log_pcsmote_x_muestra_heart_D25_R50_Pproporcion.xlsx.
# Each line: Row|F_value|G_value|K_value|Z_timestamp
$rowData = @"
2|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.955135
3|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.955135
4|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.955135
5|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.955135
6|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.955135
7|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.956135
8|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.956135
9|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.957136
10|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.957136
11|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.957136
12|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.958136
13|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.958136
14|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.958136
15|0.928947540280481|0.2857142857142857|1|2025-10-19T23:55:55.958136
16|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.959136
17|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.959136
18|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.959136
19|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.960134
20|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.960134
21|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.960134
22|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.960134
23|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.960134
24|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.973472
25|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.975479
26|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.975479
27|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.975479
28|0.928947540280481|0.2857142857142857|0.4285714285714285|2025-10-19T23:55:55.976472
29|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.976472
30|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.976472
31|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.976472
32|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.976472
33|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.977469
34|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.977469
35|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.977469
36|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.978468
37|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.978468
38|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.978468
39|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.978468
40|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.979468
41|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.979468
42|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.979468
43|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.980470
44|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.980470
45|0.928947540280481|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:55.980470
46|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.054613
47|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.055610
48|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.055610
49|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.055610
50|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.055610
51|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.055610
52|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.056610
53|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.056610
54|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.056610
55|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.056610
56|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.057612
57|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.057612
58|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.057612
59|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.057612
60|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.058611
61|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.058611
62|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.058611
63|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.058611
64|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.059610
65|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.059610
66|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.059610
67|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.059610
68|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.059610
69|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.060608
70|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.060608
71|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.060608
72|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.060608
73|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.060608
74|0.9324382446555785|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.061608
75|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.128961
76|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.128961
77|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.128961
78|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.129962
79|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.129962
80|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.129962
81|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.129962
82|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.129962
83|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.130960
84|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.130960
85|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.130960
86|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.130960
87|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.131959
88|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.133960
89|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.135963
90|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.136961
91|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.137506
92|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.137506
93|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.137506
94|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.138487
95|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.138487
96|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.138487
97|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.138487
98|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.139489
99|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.139489
100|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.139489
101|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.139489
102|0.9220334100849285|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.140490
103|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.166440
104|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.167434
105|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.167434
106|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.168431
107|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.168431
108|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.168431
109|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.169432
110|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.169432
111|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.169432
112|0.8726153446923122|0.2857142857142857|0.2857142857142857|2025-10-19T23:55:56.169432
"@

$lines = $rowData -split "`r`n|`n" | Where-Object { $_.Trim() -ne "" }

$count = 0
foreach ($line in $lines) {
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $fVal = [double]$parts[1]
    $gVal = [double]$parts[2]
    $kVal = [double]$parts[3]
    $zVal = [string]$parts[4]

    $ws.Cells.Item($row, 6).Value = $fVal
    $ws.Cells.Item($row, 7).Value = $gVal
    $ws.Cells.Item($row, 11).Value = $kVal
    $ws.Cells.Item($row, 26).Value = $zVal
    $count++
}

Write-Host "Applied updates to $count rows"
